$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" to "misc."
$ws.Name = "misc."

# New header cells K4 / L4 (same style as the rest of the row-4 header band)
$ws.Range("K4").Value = "other_indexes"
$ws.Range("K4").Style = "Heading 3"

$ws.Range("L4").Value = "commodity"
$ws.Range("L4").Style = "Heading 3"

# New row 11
$ws.Range("B11").Value = "flo_emis"
$ws.Range("D11").Value = "gas"
$ws.Range("L11").Value = "co2captured"
$ws.Range("K11").Value = "co2"
$ws.Range("E11").Value = "*ccs,*ccs-rf"
$ws.Range("H11").Value = 0.95

# New row 12
$ws.Range("B12").Value = "flo_emis"
$ws.Range("D12").Value = "coal,oil"
$ws.Range("L12").Value = "co2captured"
$ws.Range("K12").Value = "co2"
$ws.Range("E12").Value = "*ccs,*ccs-rf"
$ws.Range("H12").Value = 0.85

# column widths (engine quantizes ColumnWidth to integer pixels at MDW=6;
# choose inputs centered within the pixel bucket that reproduces the target
# OOXML widths of 10.33203125 (62px) and 12 (72px) as closely as this runtime allows)
$ws.Columns.Item(5).ColumnWidth = 9.5
$ws.Columns.Item(11).ColumnWidth = 11.17

# selection
$ws.Range("D13").Select()
